{"js": "// Merge the split \"<id>...</id>\" runs back into a single run per occurrence,\n// e.g. \"<id>\" + \"p037v_1\" + \"</id>\" (3 separate runs) -> \"<id>p037v_1</id>\" (1 run),\n// keeping the formatting (Courier New / color 7f6000 / sz 18) of the first run.\n// This mirrors the same fix for p037v_1, p037v_2 and p037v_3.\nconst ids = [\"p037v_1\", \"p037v_2\", \"p037v_3\"];\n\nfor (const id of ids) {\n  const full = \"<id>\" + id + \"</id>\";\n  const results = context.document.body.search(full, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(full, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Merge the split \"<id>...</id>\" runs back into a single run per occurrence,\n# e.g. \"<id>\" + \"p037v_1\" + \"</id>\" (3 separate runs) -> \"<id>p037v_1</id>\" (1 run),\n# keeping the formatting (Courier New / color 7f6000 / sz 18) of the first run.\n# This mirrors the same fix for p037v_1, p037v_2 and p037v_3.\n$d = $word.ActiveDocument\n$ids = @(\"p037v_1\", \"p037v_2\", \"p037v_3\")\n\nforeach ($id in $ids) {\n    $full = \"<id>\" + $id + \"</id>\"\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $full\n    $rng.Find.Replacement.Text = $full\n    $rng.Find.Execute($full, $false, $false, $false, $false, $false, $true, 1, $false, $full, 2)\n}\n"}
